$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("D2").Value = 0.52228182
$ws.Range("D3").Value = 1.0837016

$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection(1)
$series.Values = $ws.Range("D2:D8")
